$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting existing rows 7..58 down to 8..59.
$ws.Rows.Item(7).Insert(-4121)  # -4121 = xlShiftDown

# Populate the newly inserted row 7 with the new weekly data point.
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 44859
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112026
$ws.Range("G7").Value = "Haba"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 6000
$ws.Range("L7").Value = 6500
$ws.Range("M7").Value = 6250
$ws.Range("N7").Value = "`$/saco 25 kilos"
$ws.Range("O7").Value = "Provincia de Diguillín"
$ws.Range("P7").Value = 250
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
